# Generate Report for Archive
#
# 1) Update the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2, Overview!F2,
#    zh-cn!C2, de-de!C2).
# 2) Narrow the "Status" column(s) - Overview columns E & F, and the
#    "Status" column (C) on the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1) Status text change -------------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2) Column width change -------------------------------------------------
# Target stored column width is 13.4101845877511 characters; the ColumnWidth
# property (Excel "number of characters" units) that yields this value once
# quantized is 12.5.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
